$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) / Volume(1h) (E) columns with the latest scrape.
# Each target cell is forced to Text format before the write so that
# numeric-looking values (e.g. "0.2670", "220.50") keep their exact
# original formatting instead of being normalized as numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.547.31'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +5.82%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.723.60'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.11'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5359'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.12%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2670'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06611'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.75'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07743'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.629'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.721.53'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.962.52'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5853'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8328'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.99'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '27.560.84'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +5.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '220.50'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +15.15%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.733'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.52%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.097'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.25%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.86'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.728'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +14.45%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.417'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.12%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05581'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.306'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.579'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.452'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.663'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.07%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9692'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.421'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5976'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01654'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.921'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.057.03'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8544'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.55%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.50'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.868.33'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.52%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '59.08'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.254'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.95%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.007'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05252'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.11%  '
